# Updates the crypto price/volume table on Sheet1 with the latest scraped
# values (GitHub Actions data refresh). For D-column price cells whose new
# text reads as a plain number (single decimal point), force the cell to
# Text format first so Excel keeps it as a literal string (e.g. "69.70"
# instead of being normalized to the number 69.7).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range('D2').Value = '34.208.19'
$ws.Range('E2').Value = '  -0.30%  '
$ws.Range('D3').Value = '1.817.33'
$ws.Range('E3').Value = '  +1.76%  '
$ws.Range('E4').Value = '  -0.26%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '225.07'
$ws.Range('E5').Value = '  +0.23%  '
$ws.Range('E6').Value = '  +0.69%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '0.998'
$ws.Range('E7').Value = '  -0.21%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '32.11'
$ws.Range('E8').Value = '  -3.78%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.293'
$ws.Range('E9').Value = '  +4.03%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0738'
$ws.Range('E10').Value = '  +11.91%  '
$ws.Range('E11').Value = '  -0.25%  '
$ws.Range('D12').Value = '2.077.17'
$ws.Range('E12').Value = '  +1.71%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.08'
$ws.Range('E13').Value = '  -0.54%  '
$ws.Range('D14').Value = '1.809.64'
$ws.Range('E14').Value = '  +1.34%  '
$ws.Range('E15').Value = '  +1.60%  '
$ws.Range('D16').Value = '34.189.39'
$ws.Range('E16').Value = '  -0.26%  '
$ws.Range('E17').Value = '  +2.35%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.70'
$ws.Range('E18').Value = '  +1.00%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '250.47'
$ws.Range('E19').Value = '  -1.91%  '
$ws.Range('D20').Value = '0.0₃0806'
$ws.Range('E20').Value = '  +8.48%  '
$ws.Range('E21').Value = '  +6.29%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '0.996'
$ws.Range('E22').Value = '  -0.28%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.27'
$ws.Range('E23').Value = '  +2.23%  '
$ws.Range('E24').Value = '  +1.01%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '160.75'
$ws.Range('E25').Value = '  +1.87%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '16.71'
$ws.Range('E26').Value = '  +1.82%  '
$ws.Range('E27').Value = '  +3.10%  '
$ws.Range('E28').Value = '  +0.91%  '
$ws.Range('E29').Value = '  -0.23%  '
$ws.Range('E30').Value = '  +3.75%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '3.80'
$ws.Range('E31').Value = '  +0.72%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '1.22'
$ws.Range('E32').Value = '  +2.43%  '
$ws.Range('E33').Value = '  +0.66%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.90'
$ws.Range('E34').Value = '  -0.40%  '
$ws.Range('D35').Value = '1.430.65'
$ws.Range('E35').Value = '  -1.05%  '
$ws.Range('B36').Value = 'TrustWalletToken'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '1.07'
$ws.Range('E36').Value = '  +1.17%  '
$ws.Range('B37').Value = 'ImmutableX'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D37').Value = '0.644'
$ws.Range('E37').Value = '  +3.15%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '0.0190'
$ws.Range('E38').Value = '  +1.05%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '0.962'
$ws.Range('E39').Value = '  +8.18%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '81.90'
$ws.Range('E40').Value = '  -1.35%  '
$ws.Range('E41').Value = '  -3.27%  '
$ws.Range('E42').Value = '  -0.39%  '
$ws.Range('E43').Value = '  +5.15%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '6.06'
$ws.Range('E44').Value = '  +3.75%  '
$ws.Range('E45').Value = '  -1.59%  '
$ws.Range('D46').Value = '1.972.46'
$ws.Range('E46').Value = '  +1.54%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '106.28'
$ws.Range('E47').Value = '  +8.03%  '
$ws.Range('E48').Value = '  -1.40%  '
$ws.Range('B49').Value = 'InjectiveProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D49').Value = '11.97'
$ws.Range('E49').Value = '  -3.21%  '
$ws.Range('B50').Value = 'PaxDollar'
$ws.Range('C50').Value = 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp'
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '0.997'
$ws.Range('E50').Value = '  -0.32%  '
$ws.Range('E51').Value = '  +5.89%  '
